$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values look numeric,
# so Excel keeps them as exact text (preserves trailing zeros/precision)
# instead of silently converting them to floating point numbers.
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).NumberFormat = "@"

# Apply the updated cell values
$ws.Cells.Item(2, 4).Value = '57.491.46'
$ws.Cells.Item(2, 5).Value = '  -0.44%  '
$ws.Cells.Item(3, 4).Value = '3.106.07'
$ws.Cells.Item(3, 5).Value = '  +1.46%  '
$ws.Cells.Item(4, 5).Value = '  +0.01%  '
$ws.Cells.Item(5, 5).Value = '  +1.26%  '
$ws.Cells.Item(6, 4).Value = '141.31'
$ws.Cells.Item(6, 5).Value = '  -0.42%  '
$ws.Cells.Item(7, 5).Value = '  +0.02%  '
$ws.Cells.Item(8, 4).Value = '3.105.67'
$ws.Cells.Item(8, 5).Value = '  +1.47%  '
$ws.Cells.Item(9, 5).Value = '  +0.11%  '
$ws.Cells.Item(10, 5).Value = '  -0.90%  '
$ws.Cells.Item(11, 5).Value = '  +0.82%  '
$ws.Cells.Item(12, 5).Value = '  +1.51%  '
$ws.Cells.Item(13, 4).Value = '3.640.66'
$ws.Cells.Item(13, 5).Value = '  +1.46%  '
$ws.Cells.Item(14, 5).Value = '  +0.90%  '
$ws.Cells.Item(15, 5).Value = '  +0.63%  '
$ws.Cells.Item(16, 5).Value = '  +0.45%  '
$ws.Cells.Item(17, 4).Value = '57.573.77'
$ws.Cells.Item(17, 5).Value = '  -0.33%  '
$ws.Cells.Item(18, 4).Value = '3.105.89'
$ws.Cells.Item(18, 5).Value = '  +1.34%  '
$ws.Cells.Item(19, 5).Value = '  +0.32%  '
$ws.Cells.Item(20, 5).Value = '  -0.20%  '
$ws.Cells.Item(21, 5).Value = '  -0.54%  '
$ws.Cells.Item(22, 4).Value = '335.99'
$ws.Cells.Item(22, 5).Value = '  +1.45%  '
$ws.Cells.Item(23, 5).Value = '  +0.09%  '
$ws.Cells.Item(24, 5).Value = '  +2.83%  '
$ws.Cells.Item(25, 4).Value = '66.55'
$ws.Cells.Item(25, 5).Value = '  +1.33%  '
$ws.Cells.Item(26, 5).Value = '  -0.52%  '
$ws.Cells.Item(27, 5).Value = '  +0.15%  '
$ws.Cells.Item(28, 4).Value = '0.0₃0919'
$ws.Cells.Item(28, 5).Value = '  +1.72%  '
$ws.Cells.Item(29, 4).Value = '6.54'
$ws.Cells.Item(29, 5).Value = '  +2.46%  '
$ws.Cells.Item(30, 5).Value = '  +0.02%  '
$ws.Cells.Item(31, 4).Value = '7.19'
$ws.Cells.Item(31, 5).Value = '  +0.06%  '
$ws.Cells.Item(32, 5).Value = '  +2.27%  '
$ws.Cells.Item(33, 5).Value = '  +0.97%  '
$ws.Cells.Item(34, 5).Value = '  +0.36%  '
$ws.Cells.Item(35, 4).Value = '157.37'
$ws.Cells.Item(35, 5).Value = '  +1.81%  '
$ws.Cells.Item(36, 5).Value = '  +3.40%  '
$ws.Cells.Item(37, 4).Value = '6.10'
$ws.Cells.Item(37, 5).Value = '  +2.51%  '
$ws.Cells.Item(38, 4).Value = '26.99'
$ws.Cells.Item(38, 5).Value = '  -0.25%  '
$ws.Cells.Item(39, 5).Value = '  +1.05%  '
$ws.Cells.Item(40, 4).Value = '0.0661'
$ws.Cells.Item(40, 5).Value = '  -1.66%  '
$ws.Cells.Item(41, 4).Value = '3.146.83'
$ws.Cells.Item(41, 5).Value = '  +1.38%  '
$ws.Cells.Item(42, 2).Value = 'Mantle'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(42, 4).Value = '0.686'
$ws.Cells.Item(42, 5).Value = '  +4.77%  '
$ws.Cells.Item(43, 2).Value = 'Filecoin'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(43, 4).Value = '3.93'
$ws.Cells.Item(43, 5).Value = '  +0.61%  '
$ws.Cells.Item(44, 4).Value = '1.50'
$ws.Cells.Item(44, 5).Value = '  +10.78%  '
$ws.Cells.Item(45, 4).Value = '36.82'
$ws.Cells.Item(46, 4).Value = '1.00'
$ws.Cells.Item(46, 5).Value = '  +0.01%  '
$ws.Cells.Item(47, 4).Value = '2.301.04'
$ws.Cells.Item(47, 5).Value = '  +2.11%  '
$ws.Cells.Item(48, 5).Value = '  +0.69%  '
$ws.Cells.Item(49, 4).Value = '0.977'
$ws.Cells.Item(49, 5).Value = '  +4.48%  '
$ws.Cells.Item(50, 4).Value = '20.76'
$ws.Cells.Item(50, 5).Value = '  +0.87%  '
$ws.Cells.Item(51, 5).Value = '  +2.16%  '
